$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet from C_11 to C_13 (this also updates the _FilterDatabase defined name)
$ws.Name = "C_13"

# Re-apply the custom number format "#,##0.0" to the cells that use it.
# Re-setting this causes Excel to re-register the custom numFmt (freeing the
# old numFmtId and taking the lowest available custom id).
$ws.Range("C6:E6,F6,C7:E7,F7").NumberFormat = "#,##0.0"
$ws.Range("C8:F82").NumberFormat = "#,##0.0"
